# Test cases updated with wait utility and description added for all test cases.
#
# Semantic changes applied:
#   1. ResetPasswordPage!C1 test e-mail value updated from "Test@gmail.com"
#      to "test1238@gmail.com".
#   2. The ResetPasswordPage sheet becomes the active/selected tab (instead
#      of LoginPage), with its selection moved to C19.

$wb = $excel.ActiveWorkbook

$loginPage = $wb.Worksheets.Item("LoginPage")
$resetPasswordPage = $wb.Worksheets.Item("ResetPasswordPage")

# Update the test e-mail address used on the reset-password page.
$resetPasswordPage.Range("C1").Value = "test1238@gmail.com"

# Make ResetPasswordPage the active sheet/tab and move its selection to C19.
$resetPasswordPage.Activate()
$resetPasswordPage.Range("C19").Select()

# Ensure LoginPage keeps its own cell selection (C11) though it is no
# longer the active tab.
$loginPage.Range("C11").Select()
$resetPasswordPage.Activate()
